$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Model Type labels (string values)
$ws.Range("E4").Value = "gemini-1.5-pro"
$ws.Range("F4").Value = "gemini-2.0-flash-thinking-exp"
$ws.Range("J4").Value = "gemini-1.5-pro"
$ws.Range("K4").Value = "gemini-2.0-flash-thinking-exp"
$ws.Range("O4").Value = "gemini-1.5-pro"
$ws.Range("P4").Value = "gemini-2.0-flash-thinking-exp"

# Row 5 - Exact Match Mean
$ws.Range("E5").Value = 0.124
$ws.Range("F5").Value = 0.1975806451612903
$ws.Range("J5").Value = 0.104
$ws.Range("K5").Value = 0.432
$ws.Range("O5").Value = 0.176
$ws.Range("P5").Value = 0.4285714285714285

# Row 6 - F1-Score Mean
$ws.Range("E6").Value = 0.3475746031746031
$ws.Range("F6").Value = 0.3047566680631197
$ws.Range("J6").Value = 0.3539272727272727
$ws.Range("K6").Value = 0.6863746031746032
$ws.Range("O6").Value = 0.5326666666666666
$ws.Range("P6").Value = 0.5608039579468151

# Row 7 - METEOR Mean
$ws.Range("E7").Value = 0.1026602918901738
$ws.Range("F7").Value = 0.07982253108447412
$ws.Range("J7").Value = 0.1017992816128101
$ws.Range("K7").Value = 0.1361237392068599
$ws.Range("O7").Value = 0.09209420530228477
$ws.Range("P7").Value = 0.08159999192832798

# Row 8 - Bert Score Mean
$ws.Range("E8").Value = 0.5734538987874985
$ws.Range("F8").Value = 0.589771473720189
$ws.Range("J8").Value = 0.6983631743192673
$ws.Range("K8").Value = 0.8364856088161469
$ws.Range("O8").Value = 0.7167851884961128
$ws.Range("P8").Value = 0.7471670659829159
